# Update LR-pair (Atrn) sheet with new TPM-based values.
# The "Sending cluster" for the original 4 rows changes from "MuSCs" to "ECs",
# and 4 additional rows are appended for "MuSCs" as sending cluster (with its
# original ligand/receptor expression numbers), each paired against the same
# four target clusters (ECs, FAPs, MuSCs, Resolving-Mac).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics (see header row 1).
$rows = @(
    @{ Row=2; A="ECs";   B="a"; C="Atrn"; D="ECs";
       E=1; F=0.3333333333333333; G=0.045339; H=0.136017; I=0.1740293637846656; J=0.1740293637846656;
       K=3; L=1; M=11.21848733333333; N=33.655462; O=0.2956250300902271; P=0.2956250300902271;
       Q=0.508634997206; R=4.577714974854; S=0.05144743590542483; T=0.05144743590542484 },

    @{ Row=3; A="ECs";   B="a"; C="Atrn"; D="FAPs";
       E=1; F=0.3333333333333333; G=0.045339; H=0.136017; I=0.1740293637846656; J=0.1740293637846656;
       K=3; L=1; M=11.768807; N=35.306421; O=0.310126830839619; P=0.3101268308396189;
       Q=0.533585940573; R=4.802273465157; S=0.0539711750635735; T=0.0539711750635735 },

    @{ Row=4; A="ECs";   B="a"; C="Atrn"; D="MuSCs";
       E=1; F=0.3333333333333333; G=0.045339; H=0.136017; I=0.1740293637846656; J=0.1740293637846656;
       K=3; L=1; M=10.88760366666667; N=32.662811; O=0.2869057178506835; P=0.2869057178506835;
       Q=0.4936330626430001; R=4.442697563787001; S=0.04993001954373722; T=0.04993001954373723 },

    @{ Row=5; A="ECs";   B="a"; C="Atrn"; D="Resolving-Mac";
       E=1; F=0.3333333333333333; G=0.045339; H=0.136017; I=0.1740293637846656; J=0.1740293637846656;
       K=3; L=1; M=4.073469666666667; N=12.220409; O=0.1073424212194705; P=0.1073424212194704;
       Q=0.184687041217; R=1.662183370953; S=0.01868073327193003; T=0.01868073327193003 },

    @{ Row=6; A="MuSCs"; B="a"; C="Atrn"; D="ECs";
       E=1; F=0.3333333333333333; G=0.215186; H=0.645558; I=0.8259706362153344; J=0.8259706362153345;
       K=3; L=1; M=11.21848733333333; N=33.655462; O=0.2956250300902271; P=0.2956250300902271;
       Q=2.414061415310667; R=21.726552737796; S=0.2441775941848022; T=0.2441775941848023 },

    @{ Row=7; A="MuSCs"; B="a"; C="Atrn"; D="FAPs";
       E=1; F=0.3333333333333333; G=0.215186; H=0.645558; I=0.8259706362153344; J=0.8259706362153345;
       K=3; L=1; M=11.768807; N=35.306421; O=0.310126830839619; P=0.3101268308396189;
       Q=2.532482503102; R=22.792342527918; S=0.2561556557760455; T=0.2561556557760455 },

    @{ Row=8; A="MuSCs"; B="a"; C="Atrn"; D="MuSCs";
       E=1; F=0.3333333333333333; G=0.215186; H=0.645558; I=0.8259706362153344; J=0.8259706362153345;
       K=3; L=1; M=10.88760366666667; N=32.662811; O=0.2869057178506835; P=0.2869057178506835;
       Q=2.342859882615334; R=21.085738943538; S=0.2369756983069463; T=0.2369756983069463 },

    @{ Row=9; A="MuSCs"; B="a"; C="Atrn"; D="Resolving-Mac";
       E=1; F=0.3333333333333333; G=0.215186; H=0.645558; I=0.8259706362153344; J=0.8259706362153345;
       K=3; L=1; M=4.073469666666667; N=12.220409; O=0.1073424212194705; P=0.1073424212194704;
       Q=0.8765536436913334; R=7.888982793222; S=0.08866168794754042; T=0.08866168794754042 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $r.A
    $ws.Cells.Item($row, 2).Value2  = $r.B
    $ws.Cells.Item($row, 3).Value2  = $r.C
    $ws.Cells.Item($row, 4).Value2  = $r.D
    $ws.Cells.Item($row, 5).Value2  = $r.E
    $ws.Cells.Item($row, 6).Value2  = $r.F
    $ws.Cells.Item($row, 7).Value2  = $r.G
    $ws.Cells.Item($row, 8).Value2  = $r.H
    $ws.Cells.Item($row, 9).Value2  = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $r.T
}
